function ConvertTo-PointsForExactEmu {
    param([double]$TargetEmu)

    # PowerPoint's Shape.Left/Top/Width/Height are expressed in points, and the
    # host stores them with single-precision (float32) rounding internally.
    # A naive $TargetEmu/12700 conversion can therefore land the re-serialized
    # EMU one unit away from the intended value after the float32 round trip.
    # Search nearby point values for one whose float32 round trip reproduces
    # the exact target EMU value.
    $basePt = $TargetEmu / 12700.0
    $f = [float]$basePt
    $emuOut = [int64]([double]$f * 12700.0)
    if ($emuOut -eq $TargetEmu) {
        return $basePt
    }

    $scales = @(1, 2, 5, 10, 20, 50, 100, 200, 500, 1000, 2000, 5000, 10000)
    foreach ($scale in $scales) {
        for ($i = -2000; $i -le 2000; $i++) {
            $cand = $basePt + ($i / 12700.0 / $scale)
            $f = [float]$cand
            $emuOut = [int64]([double]$f * 12700.0)
            if ($emuOut -eq $TargetEmu) {
                return $cand
            }
        }
    }
    return $basePt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

$shp.Left   = ConvertTo-PointsForExactEmu(-50)
$shp.Top    = ConvertTo-PointsForExactEmu(7)
$shp.Width  = ConvertTo-PointsForExactEmu(13817699)
$shp.Height = ConvertTo-PointsForExactEmu(10677299)
